$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add four new "Display CF All Data" verification rows (21-24) underneath the
# existing EmulatorData test matrix (rows 1-20). These are new Katalon test
# fixtures exercising the Convenience-Fee display logic at two MV tiers
# (2.3 / 2.5) crossed with two Amount thresholds (5 / 200).
# ---------------------------------------------------------------------------

# Seed the brand-new lookup values first (this is also the order in which
# they were originally authored) so the shared-string table grows the same
# way it did in the source edit.
$ws.Range("A21").Value = "Display CF All Data"
$ws.Range("C21").Value = "20"
$ws.Range("C22").Value = "21"
$ws.Range("F21").Value = "26413419"
$ws.Range("F23").Value = "26413420"
$ws.Range("C23").Value = "22"
$ws.Range("E22").Value = "200"
$ws.Range("C24").Value = "23"

function Set-EmulatorRow($Sheet, $Row, $MV, $Amount, $Can) {
    $Sheet.Range("A$Row").Value = "Display CF All Data"
    $Sheet.Range("D$Row").Value = $MV
    $Sheet.Range("E$Row").Value = $Amount
    $Sheet.Range("F$Row").Value = $Can
    $Sheet.Range("G$Row").Value = "PayNow"
    $Sheet.Range("H$Row").Value = "en_US"
    $Sheet.Range("I$Row").Value = "Elizath"
    $Sheet.Range("J$Row").Value = "Christine"
    $Sheet.Range("K$Row").Value = "258 Underwood rd"
    $Sheet.Range("L$Row").Value = "Suite 600"
    $Sheet.Range("M$Row").Value = "840"
    $Sheet.Range("N$Row").Value = "Arlington"
    $Sheet.Range("O$Row").Value = "VA"
    $Sheet.Range("P$Row").Value = "22201"
    $Sheet.Range("R$Row").Value = "Some Company"
    $Sheet.Range("S$Row").Value = "iahmed@govolution.com"
    $Sheet.Range("T$Row").Value = "udf data 1"
    $Sheet.Range("U$Row").Value = "udf data 2"
    $Sheet.Range("V$Row").Value = "udf data 3"
    $Sheet.Range("W$Row").Value = "udf data 4"
    $Sheet.Range("X$Row").Value = "udf data 5"
    $Sheet.Range("Y$Row").Value = "udf data 6"
    $Sheet.Range("AB$Row").Value = "udf data 9"
    $Sheet.Range("AC$Row").Value = "udf data 10"

    # Column S carries the bordered/unwrapped "email" cell style used
    # throughout the sheet (style index 2) instead of the default text
    # style (index 1) that a bare value-assignment would pick up.
    $Sheet.Range("S2").Copy()
    $Sheet.Range("S$Row").PasteSpecial(-4122)
}

Set-EmulatorRow $ws 21 "2.3" "5"   "26413419"
Set-EmulatorRow $ws 22 "2.3" "200" "26413419"
Set-EmulatorRow $ws 23 "2.5" "5"   "26413420"
Set-EmulatorRow $ws 24 "2.5" "200" "26413420"

$excel.CutCopyMode = $false

# Move the active selection below the newly-added rows, matching the saved
# worksheet view state of the edited workbook.
$ws.Range("A25:AG25").Select()
